$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows ---
$ws.Range("H273").Value = 25977

$ws.Range("H279").Value = 42220

$ws.Range("H280").Value = 32670
$ws.Range("I280").Value = 2130

$ws.Range("H281").Value = 41465
$ws.Range("I281").Value = 2940

# --- New rows 282-284 ---
$data = @(
    @(44176, 130794, 95416, 34230, 17956, 3707, 1148, 42243, 2545),
    @(44177, 132984, 97068, 34741, 9991, 2190, 1175, 15304, 911),
    @(44178, 133489, 98585, 33699, 2899, 505, 1205, 1509, 118)
)

$row = 282
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 1).NumberFormat = "yyyy-mm-dd"
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $ws.Cells.Item($row, 6).Value = $entry[5]
    $ws.Cells.Item($row, 7).Value = $entry[6]
    $ws.Cells.Item($row, 8).Value = $entry[7]
    $ws.Cells.Item($row, 9).Value = $entry[8]
    $row++
}
